$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Part 1: remove the "Requirements for SDS Diffusion phase:" paragraph and its
# three bullet points, and relocate the (hidden) _GoBack bookmark from the
# blank paragraph after "Traversal is Depth First Search" onto the
# "Requirements for Agents" paragraph that now immediately follows the
# deleted block.
# ---------------------------------------------------------------------------

$found = $d.Content.Find.Execute("Requirements for SDS Diffusion phase:")
if (-not $found) {
    throw "Could not find 'Requirements for SDS Diffusion phase:' paragraph"
}
$startPara = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.Start).Paragraphs.Item(1)

# Locate the paragraph boundaries for the heading + the 3 bullets that follow it.
$deleteStart = -1
$deleteEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Requirements for SDS Diffusion phase:") {
        $deleteStart = $p.Range.Start
        $deleteEnd = $d.Paragraphs.Item($i + 3).Range.End
        break
    }
}
if ($deleteStart -lt 0) {
    throw "Could not locate paragraph range to delete"
}
$d.Range($deleteStart, $deleteEnd).Delete()

# Find the now-adjacent "Requirements for Agents" paragraph and the paragraph
# that still carries the _GoBack bookmark (the blank ListParagraph right after
# "Traversal is Depth First Search").
$targetIdx = -1
$bookmarkIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "^Requirements for Agents") {
        $targetIdx = $i
    }
    if ($t -match "^Traversal is Depth First Search") {
        $bookmarkIdx = $i + 1
    }
}
if ($targetIdx -lt 0 -or $bookmarkIdx -lt 0) {
    throw "Could not locate target / bookmark paragraphs"
}

$bookmarkPara = $d.Paragraphs.Item($bookmarkIdx)
$xmlNoBookmark = "<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
$bookmarkPara.Range.InsertXML($xmlNoBookmark)

$targetPara = $d.Paragraphs.Item($targetIdx)
$xmlWithBookmark = "<w:p $ns><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/><w:r><w:t>Requirements for Agents</w:t></w:r></w:p>"
$targetPara.Range.InsertXML($xmlWithBookmark)

# ---------------------------------------------------------------------------
# Part 2: the "Error with test for random hypothesis..." paragraph currently
# opens with a lastRenderedPageBreak marker. Move that marker so it instead
# opens the following "Error when checking..." paragraph.
# ---------------------------------------------------------------------------

$agentTestIdx = -1
$checkingIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "^Error with test for random") {
        $agentTestIdx = $i
    }
    if ($t -match "^Error when checking if a single node graph") {
        $checkingIdx = $i
    }
}
if ($agentTestIdx -lt 0 -or $checkingIdx -lt 0) {
    throw "Could not locate the two 'Error ...' paragraphs"
}

$pAgentTest = $d.Paragraphs.Item($agentTestIdx)
$xmlAgentTest = "<w:p $ns>" + `
  "<w:r><w:t xml:space='preserve'>Error with test for random </w:t></w:r>" + `
  "<w:proofErr w:type='gramStart'/>" + `
  "<w:r><w:t>hypothesis,</w:t></w:r>" + `
  "<w:proofErr w:type='gramEnd'/>" + `
  "<w:r><w:t xml:space='preserve'> was checking class EWG against class EWG which resulted in different EWG which had the same nodes and edges. </w:t></w:r>" + `
  "<w:proofErr w:type='gramStart'/>" + `
  "<w:r><w:t>Changed test to be based on the weight of the graphs.</w:t></w:r>" + `
  "<w:proofErr w:type='gramEnd'/>" + `
  "<w:r><w:t xml:space='preserve'> (Agent test)</w:t></w:r>" + `
  "</w:p>"
$pAgentTest.Range.InsertXML($xmlAgentTest)

$pChecking = $d.Paragraphs.Item($checkingIdx)
$xmlChecking = "<w:p $ns>" + `
  "<w:r><w:lastRenderedPageBreak/><w:t>Error when checking if a single node graph with no edges is a spanning tree, the result returned was true when it is supposed to be false</w:t></w:r>" + `
  "<w:r><w:t>. Revised the if statement that set spanning tree check to true and added a condition for single node graphs as they will not be spanning trees or contain cycles due to</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> the</w:t></w:r>" + `
  "<w:r><w:t xml:space='preserve'> graph not allowing for self-looping vertices.</w:t></w:r>" + `
  "</w:p>"
$pChecking.Range.InsertXML($xmlChecking)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
